$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23.336.55"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.632.29"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.66%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("E5").Value = "  -0.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "302.97"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.20%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3816"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.82%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "51.98"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.33%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3552"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08126"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.73%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.222"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.16%  "
$ws.Range("E12").Value = "  +0.02%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.24"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.30%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.425"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.69%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.290"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.30%  "
$ws.Range("E16").Value = "  -1.48%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.626.92"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.31%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.68"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.28%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06931"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.41%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.575"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.67%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.30"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.10%  "
$ws.Range("E22").Value = "  -0.14%  "
$ws.Range("E23").Value = "  -2.28%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "23.332.79"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.554"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.43%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.127"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.30%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.25%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "151.54"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.16%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.267"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.36%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "132.99"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.15%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.808.51"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.48%  "
$ws.Range("B32").Value = "WEMIXTOKEN"
$ws.Range("C32").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.146"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.41%  "
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.073"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +13.50%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.499"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.72%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "11.55"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.96%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02736"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2486"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.08732"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.94%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.924"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.59%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06956"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.00%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6941"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.85%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.318"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.75%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.12"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.37"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.76%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.000"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6381"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.37%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.267"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.51%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.955"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.76%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07923"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.66%  "
$ws.Range("E50").Value = "  +3.57%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.181"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.94%  "
